$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.64%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.37%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.769"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08113"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.24%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.097"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.60%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.730"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.521"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.19%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.954"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.50%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9243"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.84%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1284"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.80%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1957"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.842"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "14.93%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09178"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.61%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03758"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "8.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1054"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.32%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001300"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.16%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006323"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.51%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.379"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.10%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3503"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.94%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1381"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.12%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2604"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.51%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.27%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004492"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.94%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.22%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02795"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.02%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05540"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.57%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007642"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.09%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009895"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.39%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1423"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.40%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002223"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.76%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01185"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.84%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006790"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.09%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.23%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002993"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "4.33%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002281"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "26.56%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.23%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.23%"
